# Cập nhật báo cáo: bổ sung tiền tố "Phần mềm" cho tên phần mềm
# trong bảng tham khảo (cột B, các dòng "Nén mp3" và "Thu âm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Phần mềm EKOS MP3Minimizer"
$ws.Range("B4").Value = "Phần mềm Mp3 Audio Editor"
